$d = $word.ActiveDocument

# Locate the bullet paragraph "Chức năng thông báo người đang chờ bổ sung
# thông tin cho người quản trị." - the new bullet needs to be inserted
# right after it, inheriting the same list formatting (ListParagraph style
# + numPr numbering).
$rng = $d.Content
$found = $rng.Find.Execute(
    "Chức năng thông báo người đang chờ bổ sung thông tin cho người quản trị.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Work out the 1-based paragraph index of the match (before mutating the
    # document) so we can reliably grab the freshly-inserted paragraph below.
    $paraIndex = $d.Range(0, $rng.Start).Paragraphs.Count + 1

    # Collapse to the end of the matched text and insert a brand-new,
    # empty paragraph right after it; Word copies the paragraph's list
    # formatting (pStyle=ListParagraph, numPr ilvl/numId) onto it.
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()

    # Fetch the newly created paragraph (the one right after the original)
    # and fill in its text.
    $newPara = $d.Paragraphs.Item($paraIndex + 1)
    $newPara.Range.InsertBefore("Thêm chức năng xác thực client để kết nối bắt đầu quét thẻ điểm danh. (Chống truy cập không được phép từ bên ngoài và hạn chế tấn công dos)")
}
